$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.001.32"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.211.85"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'602.41"
$ws.Range("E5").Value = "  +3.91%  "
$ws.Range("D6").Value = "'153.89"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.210.23"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  +1.80%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "'39.50"
$ws.Range("E14").Value = "  +5.13%  "
$ws.Range("D15").Value = "3.739.27"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "'7.51"
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("D17").Value = "66.082.31"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "3.212.47"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'511.35"
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("D22").Value = "'0.740"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "'15.57"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'8.10"
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").Value = "'85.11"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'9.37"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  +4.43%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = "  +7.59%  "
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'55.23"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "'485.25"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("D39").Value = "'0.0421"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "'0.305"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("D44").Value = "2.953.20"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "0.0₃0640"
$ws.Range("E46").Value = "  +5.23%  "
$ws.Range("D47").Value = "'28.73"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").Value = "'120.29"
$ws.Range("E51").Value = "  +0.08%  "
